$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.7476426666666667
$ws.Range("H2").Value = 2.242928
$ws.Range("I2").Value = 0.001581772089386036
$ws.Range("J2").Value = 0.001581772089386036
$ws.Range("M2").Value = 3.303267
$ws.Range("N2").Value = 9.909801000000002
$ws.Range("O2").Value = 0.03362563178859915
$ws.Range("P2").Value = 0.03362563178859915
$ws.Range("Q2").Value = 2.469663348592
$ws.Range("R2").Value = 22.226970137328
$ws.Range("S2").Value = 0.00005318808585117799
$ws.Range("T2").Value = 0.000053188085851178

$ws.Range("G3").Value = 0.7476426666666667
$ws.Range("H3").Value = 2.242928
$ws.Range("I3").Value = 0.001581772089386036
$ws.Range("J3").Value = 0.001581772089386036
$ws.Range("M3").Value = 37.82684066666667
$ws.Range("O3").Value = 0.3850586149964086
$ws.Range("P3").Value = 0.3850586149964086
$ws.Range("Q3").Value = 28.28096002760178
$ws.Range("R3").Value = 254.528640248416
$ws.Range("S3").Value = 0.0006090749699789624
$ws.Range("T3").Value = 0.0006090749699789625

$ws.Range("G4").Value = 0.7476426666666667
$ws.Range("H4").Value = 2.242928
$ws.Range("I4").Value = 0.001581772089386036
$ws.Range("J4").Value = 0.001581772089386036
$ws.Range("M4").Value = 9.149396
$ws.Range("N4").Value = 27.448188
$ws.Range("O4").Value = 0.09313634682999644
$ws.Range("P4").Value = 0.09313634682999644
$ws.Range("Q4").Value = 6.840478823829333
$ws.Range("R4").Value = 61.564309414464
$ws.Range("S4").Value = 0.000147320473923066
$ws.Range("T4").Value = 0.000147320473923066

$ws.Range("G5").Value = 0.7476426666666667
$ws.Range("H5").Value = 2.242928
$ws.Range("I5").Value = 0.001581772089386036
$ws.Range("J5").Value = 0.001581772089386036
$ws.Range("M5").Value = 47.95707433333333
$ws.Range("N5").Value = 143.871223
$ws.Range("O5").Value = 0.4881794063849957
$ws.Range("P5").Value = 0.4881794063849957
$ws.Range("Q5").Value = 35.85475494010489
$ws.Range("R5").Value = 322.692794460944
$ws.Range("S5").Value = 0.0007721885596328293
$ws.Range("T5").Value = 0.0007721885596328294

$ws.Range("I6").Value = 0.02590993131491687
$ws.Range("J6").Value = 0.02590993131491688
$ws.Range("M6").Value = 3.303267
$ws.Range("N6").Value = 9.909801000000002
$ws.Range("O6").Value = 0.03362563178859915
$ws.Range("P6").Value = 0.03362563178859915
$ws.Range("Q6").Value = 40.45387332496401
$ws.Range("R6").Value = 364.0848599246761
$ws.Range("S6").Value = 0.0008712378100632895
$ws.Range("T6").Value = 0.0008712378100632896

$ws.Range("I7").Value = 0.02590993131491687
$ws.Range("J7").Value = 0.02590993131491688
$ws.Range("M7").Value = 37.82684066666667
$ws.Range("O7").Value = 0.3850586149964086
$ws.Range("P7").Value = 0.3850586149964086
$ws.Range("Q7").Value = 463.2511451883636
$ws.Range("S7").Value = 0.009976842266773969
$ws.Range("T7").Value = 0.00997684226677397

$ws.Range("I8").Value = 0.02590993131491687
$ws.Range("J8").Value = 0.02590993131491688
$ws.Range("M8").Value = 9.149396
$ws.Range("N8").Value = 27.448188
$ws.Range("O8").Value = 0.09313634682999644
$ws.Range("P8").Value = 0.09313634682999644
$ws.Range("Q8").Value = 112.0492248382987
$ws.Range("R8").Value = 1008.443023544688
$ws.Range("S8").Value = 0.002413156349287484
$ws.Range("T8").Value = 0.002413156349287484

$ws.Range("I9").Value = 0.02590993131491687
$ws.Range("J9").Value = 0.02590993131491688
$ws.Range("M9").Value = 47.95707433333333
$ws.Range("N9").Value = 143.871223
$ws.Range("O9").Value = 0.4881794063849957
$ws.Range("P9").Value = 0.4881794063849957
$ws.Range("Q9").Value = 587.3123214431498
$ws.Range("R9").Value = 5285.810892988347
$ws.Range("S9").Value = 0.01264869488879213
$ws.Range("T9").Value = 0.01264869488879213

$ws.Range("G10").Value = 18.93023433333333
$ws.Range("H10").Value = 56.79070299999999
$ws.Range("I10").Value = 0.04005030430848061
$ws.Range("J10").Value = 0.04005030430848062
$ws.Range("M10").Value = 3.303267
$ws.Range("N10").Value = 9.909801000000002
$ws.Range("O10").Value = 0.03362563178859915
$ws.Range("P10").Value = 0.03362563178859915
$ws.Range("Q10").Value = 62.531618375567
$ws.Range("R10").Value = 562.784565380103
$ws.Range("S10").Value = 0.001346716785698315
$ws.Range("T10").Value = 0.001346716785698316

$ws.Range("G11").Value = 18.93023433333333
$ws.Range("H11").Value = 56.79070299999999
$ws.Range("I11").Value = 0.04005030430848061
$ws.Range("J11").Value = 0.04005030430848062
$ws.Range("M11").Value = 37.82684066666667
$ws.Range("O11").Value = 0.3850586149964086
$ws.Range("P11").Value = 0.3850586149964086
$ws.Range("Q11").Value = 716.0709579096629
$ws.Range("R11").Value = 6444.638621186966
$ws.Range("S11").Value = 0.01542171470720824
$ws.Range("T11").Value = 0.01542171470720824

$ws.Range("G12").Value = 18.93023433333333
$ws.Range("H12").Value = 56.79070299999999
$ws.Range("I12").Value = 0.04005030430848061
$ws.Range("J12").Value = 0.04005030430848062
$ws.Range("M12").Value = 9.149396
$ws.Range("N12").Value = 27.448188
$ws.Range("O12").Value = 0.09313634682999644
$ws.Range("P12").Value = 0.09313634682999644
$ws.Range("Q12").Value = 173.2002102884626
$ws.Range("R12").Value = 1558.801892596164
$ws.Range("S12").Value = 0.003730139032721551
$ws.Range("T12").Value = 0.003730139032721552

$ws.Range("G13").Value = 18.93023433333333
$ws.Range("H13").Value = 56.79070299999999
$ws.Range("I13").Value = 0.04005030430848061
$ws.Range("J13").Value = 0.04005030430848062
$ws.Range("M13").Value = 47.95707433333333
$ws.Range("N13").Value = 143.871223
$ws.Range("O13").Value = 0.4881794063849957
$ws.Range("P13").Value = 0.4881794063849957
$ws.Range("Q13").Value = 907.8386550710853
$ws.Range("R13").Value = 8170.547895639767
$ws.Range("S13").Value = 0.0195517337828525
$ws.Range("T13").Value = 0.0195517337828525

$ws.Range("G14").Value = 440.7369333333333
$ws.Range("H14").Value = 1322.2108
$ws.Range("I14").Value = 0.9324579922872165
$ws.Range("J14").Value = 0.9324579922872166
$ws.Range("M14").Value = 3.303267
$ws.Range("N14").Value = 9.909801000000002
$ws.Range("O14").Value = 0.03362563178859915
$ws.Range("P14").Value = 0.03362563178859915
$ws.Range("Q14").Value = 1455.8717675612
$ws.Range("R14").Value = 13102.8459080508
$ws.Range("S14").Value = 0.03135448910698637
$ws.Range("T14").Value = 0.03135448910698638

$ws.Range("G15").Value = 440.7369333333333
$ws.Range("H15").Value = 1322.2108
$ws.Range("I15").Value = 0.9324579922872165
$ws.Range("J15").Value = 0.9324579922872166
$ws.Range("M15").Value = 37.82684066666667
$ws.Range("O15").Value = 0.3850586149964086
$ws.Range("P15").Value = 0.3850586149964086
$ws.Range("Q15").Value = 16671.68575311529
$ws.Range("R15").Value = 150045.1717780376
$ws.Range("S15").Value = 0.3590509830524475
$ws.Range("T15").Value = 0.3590509830524475

$ws.Range("G16").Value = 440.7369333333333
$ws.Range("H16").Value = 1322.2108
$ws.Range("I16").Value = 0.9324579922872165
$ws.Range("J16").Value = 0.9324579922872166
$ws.Range("M16").Value = 9.149396
$ws.Range("N16").Value = 27.448188
$ws.Range("O16").Value = 0.09313634682999644
$ws.Range("P16").Value = 0.09313634682999644
$ws.Range("Q16").Value = 4032.476734892266
$ws.Range("R16").Value = 36292.2906140304
$ws.Range("S16").Value = 0.08684573097406434
$ws.Range("T16").Value = 0.08684573097406435

$ws.Range("G17").Value = 440.7369333333333
$ws.Range("H17").Value = 1322.2108
$ws.Range("I17").Value = 0.9324579922872165
$ws.Range("J17").Value = 0.9324579922872166
$ws.Range("M17").Value = 47.95707433333333
$ws.Range("N17").Value = 143.871223
$ws.Range("O17").Value = 0.4881794063849957
$ws.Range("P17").Value = 0.4881794063849957
$ws.Range("Q17").Value = 21136.45387331204
$ws.Range("R17").Value = 190228.0848598084
$ws.Range("S17").Value = 0.4552067891537182
$ws.Range("T17").Value = 0.4552067891537183

